$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141; this shifts the existing rows 141-196
# down to 142-197 (content + formatting moves with them automatically).
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new record.
$ws.Cells.Item(141, 1).Value = 4
$ws.Cells.Item(141, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(141, 3).Value = "Los Lagos"
$ws.Cells.Item(141, 4).Value = 44553
$ws.Cells.Item(141, 5).Value = 10
$ws.Cells.Item(141, 6).Value = "Fruta"
$ws.Cells.Item(141, 7).Value = 100104
$ws.Cells.Item(141, 8).Value = "Frutos de pepita"
$ws.Cells.Item(141, 9).Value = 100104005
$ws.Cells.Item(141, 10).Value = "Pera"
$ws.Cells.Item(141, 11).Value = "Packham's Triumph"
$ws.Cells.Item(141, 12).Value = "Primera"
$ws.Cells.Item(141, 13).Value = 600
$ws.Cells.Item(141, 14).Value = 13000
$ws.Cells.Item(141, 15).Value = 14000
$ws.Cells.Item(141, 16).Value = 13500
$ws.Cells.Item(141, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(141, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(141, 19).Value = 900
$ws.Cells.Item(141, 20).Value = 15
